$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update selection: active cell A2, selected range A2:L130 ---
$ws.Range("A2").Activate() | Out-Null
$ws.Range("A2:L130").Select() | Out-Null

# --- Update the financed-rate (H) table and its dependent columns (J, K, L)
# for rows 57-62 (MiPyme / No MiPyme, MACRO, Tasa 0%) and rows 121-125
# (AGRONACIÓN, Tarjeta, Financia Banco). Columns J/K/L are plain cached
# numbers (no formulas in this sheet), so each one is written explicitly.

$rows = @{
    57  = @{ H = 0.42499999999999999; J = 20.958904109589042; K = 134.14342465753427; L = 0.69235277777777804 }
    58  = @{ H = 0.49;                 J = 36.246575342465754; K = 151.09745205479453; L = 0.6907618518518519 }
    59  = @{ H = 0.53;                 J = 52.273972602739725; K = 168.87183561643837; L = 0.69828388888888904 }
    60  = @{ H = 0.42499999999999999; J = 20.958904109589042; K = 134.14342465753427; L = 0.69235277777777804 }
    61  = @{ H = 0.49;                 J = 36.246575342465754; K = 151.09745205479453; L = 0.6907618518518519 }
    62  = @{ H = 0.53;                 J = 52.273972602739725; K = 168.87183561643837; L = 0.69828388888888904 }
    121 = @{ H = 0.42; J = 6.9041095890410951;  K = 109.04219178082191; L = 0.55006666666666593 }
    122 = @{ H = 0.42; J = 10.356164383561643;  K = 112.56328767123287; L = 0.50951111111111091 }
    123 = @{ H = 0.42; J = 20.712328767123285;  K = 123.12657534246574; L = 0.46895555555555524 }
    124 = @{ H = 0.42; J = 31.06849315068493;   K = 133.68986301369861; L = 0.45543703703703675 }
    125 = @{ H = 0.42; J = 41.42465753424657;   K = 144.25315068493151; L = 0.4486777777777779 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 8).Value = $vals.H
    $ws.Cells.Item($r, 10).Value = $vals.J
    $ws.Cells.Item($r, 11).Value = $vals.K
    $ws.Cells.Item($r, 12).Value = $vals.L
}
